$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a date-looking string into column A without letting Excel's
# autoconvert turn it into a serial date number, and without leaving an
# explicit (non-default) cell style behind.
function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Update existing rows 86 and 87 (values only, no new style needed)
$ws.Range("B86").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE']"
$ws.Range("B87").Value = "['BTCUSD.SPOT']"

# Append new rows 88-95
Set-TextCell "A88" "2025-09-07"
$ws.Range("B88").Value = "['BTCUSD.SPOT']"

Set-TextCell "A89" "2025-09-08"
$ws.Range("B89").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE']"

Set-TextCell "A90" "2025-09-09"
$ws.Range("B90").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE']"

Set-TextCell "A91" "2025-09-10"
$ws.Range("B91").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE']"

Set-TextCell "A92" "2025-09-11"
$ws.Range("B92").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.VOLSURFACE']"

Set-TextCell "A93" "2025-09-12"
$ws.Range("B93").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT']"

Set-TextCell "A94" "2025-09-13"
$ws.Range("B94").Value = "[]"

Set-TextCell "A95" "2025-09-14"
$ws.Range("B95").Value = "[]"
